$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with P1=14, Q1=15, matching the style of O1 (bold/centered/bordered)
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New values for columns B..Q, rows 2..25 (one row per data point, B..K recomputed,
# L..N stay 0, O becomes 0, P is a new all-zero column, Q is a new results column)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

$rowsData = @(
  @(3.576352370229188, 1.000931929389992, 0.06682399138679784, 1.399864355909997, 0.2932408734036258, 0.1904135290156077, 0.01639006590631475, 0, 0.1619758841354155, 0.0635532518251889, 0, 0, 0, 0, 0, 0.6941096359857539),
  @(3.12248155556523, 0.8795569721254708, 0.05855141369686834, 1.223823267901395, 0.268544262242294, 0.173926500414872, 0.01256638071878891, 0, 0.1597059478341478, 0.0722535752831659, 0, 0, 0, 0, 0, 0.6540692518301512),
  @(2.84348365213458, 0.8049642780964916, 0.05345631561861808, 1.115951810711422, 0.254110875639121, 0.1644552365826826, 0.01040243372296885, 0, 0.1587267520974436, 0.07813653803829368, 0, 0, 0, 0, 0, 0.6317260372592557),
  @(2.729687759135516, 0.7745448471711427, 0.05137584598264766, 1.072028594873473, 0.2484002147420128, 0.1607474402720825, 0.009562658541953366, 0, 0.1584265072349922, 0.08066284490936848, 0, 0, 0, 0, 0, 0.6231476329705714),
  @(2.710785279339746, 0.7694922332246961, 0.05103012721533418, 1.064736751419957, 0.2474619369075413, 0.1601405894707568, 0.009425635418719512, 0, 0.1583824630524262, 0.08108992346807398, 0, 0, 0, 0, 0, 0.6217539493999311),
  @(2.841949401205511, 0.8045541278775374, 0.05342827486794732, 1.115359329361169, 0.2540331833275644, 0.1644046328288198, 0.01039094335275727, 0, 0.1587223099534825, 0.07817009570252687, 0, 0, 0, 0, 0, 0.6316082613942484),
  @(3.419910274594997, 0.9590926440747296, 0.06397471847938618, 1.339107884202264, 0.284566116672309, 0.1845860700272581, 0.01503094827811002, 0, 0.1611038105262068, 0.06643650949812052, 0, 0, 0, 0, 0, 0.6798156205158961),
  @(4.551878420986441, 1.26187530566267, 0.08454422077662116, 1.78056183305948, 0.3508347246181103, 0.2299065146863057, 0.02578465778696426, 0, 0.1693213920722201, 0.04805135562994423, 0, 0, 0, 0, 0, 0.7938924287182942),
  @(5.383279404422581, 1.480951481558975, 0.1008814953904533, 2.004452670820953, 0.3977106019286154, 0.2612575997408015, 0.0341479099069355, 0, 0.1748412217086397, 0.03685977187935929, 0, 0, 0, 0, 0, 0.872447595844136),
  @(5.752464130054989, 1.549011572455925, 0.1199402269433847, 1.307533055275542, 0.3641612700215759, 0.2241168566502978, 0.04909507406784996, 0, 0.1522468126925673, 0.02569406017226683, 0, 0, 0, 0, 0, 0.7465609327796159),
  @(5.888975715355514, 1.56001835655303, 0.1333584267596706, 0.8049952184876332, 0.3270813402513113, 0.188665350680175, 0.08455913375654234, 0, 0.1329010134908728, 0.02060106516065696, 0, 0, 0, 0, 0, 0.6309332823787486),
  @(5.852980200396416, 1.527638864193818, 0.1430535465160716, 0.4230316865514183, 0.2850941363040249, 0.1524124110741241, 0.1371064354420639, 0, 0.1148509520155514, 0.01953425138549778, 0, 0, 0, 0, 0, 0.5162253271228394),
  @(5.749771716033081, 1.485498708855062, 0.1481590256553034, 0.2287956058398137, 0.2539694725339885, 0.1272837651619199, 0.1843406214986487, 0, 0.1030423497333075, 0.02083997426219009, 0, 0, 0, 0, 0, 0.4381444376191581),
  @(5.687754969891103, 1.466054100094141, 0.1484794892514998, 0.1905170396152158, 0.2452431033458069, 0.1207019738358071, 0.1961392119646348, 0, 0.1003370084451589, 0.02172219701502209, 0, 0, 0, 0, 0, 0.4184843424691849),
  @(5.336158623809467, 1.378923685387065, 0.1390637338417804, 0.1865093761526744, 0.2338899834970363, 0.1148450827808745, 0.1807371326037952, 0, 0.1015433602594271, 0.02582143384408075, 0, 0, 0, 0, 0, 0.4083773421332637),
  @(5.121953701495841, 1.334502587289137, 0.1291561061871676, 0.2728430739711172, 0.2416319639306082, 0.1235016168484506, 0.1419553923922194, 0, 0.1087172156465144, 0.02832280197602444, 0, 0, 0, 0, 0, 0.4416274316102147),
  @(5.001188400532328, 1.321652906777899, 0.1179930550577808, 0.5104868431844807, 0.2672687853399083, 0.1468541838353019, 0.08995275647465917, 0, 0.1222810021007348, 0.0305663628150139, 0, 0, 0, 0, 0, 0.5194862213739668),
  @(4.964501139649542, 1.335299708563696, 0.1070144088365907, 0.9512415835460075, 0.306886307043996, 0.1826687615265996, 0.04783947179917192, 0, 0.140896115214467, 0.0338536059382859, 0, 0, 0, 0, 0, 0.6344116587011115),
  @(5.164812890301505, 1.423254119543174, 0.09664310183642044, 1.941202215490463, 0.3847791124496567, 0.252444414766515, 0.03180764788598989, 0, 0.1730787705045884, 0.03950605051602896, 0, 0, 0, 0, 0, 0.8499353941120091),
  @(5.805335792708206, 1.597189942350894, 0.1072058252821364, 2.274308694021386, 0.4326995008928378, 0.287844598568455, 0.04001844929238718, 0, 0.1829615619393365, 0.03364729787374454, 0, 0, 0, 0, 0, 0.9460889171629674),
  @(6.224150600814596, 1.709217398090004, 0.1147530322703147, 2.440414567681472, 0.4618342928692414, 0.3088576357552171, 0.04528532151681874, 0, 0.1884540694698131, 0.03007433856349095, 0, 0, 0, 0, 0, 1.002394635086119),
  @(6.000537257581982, 1.649404860288314, 0.1107249151151564, 2.351657497299655, 0.4461723037852678, 0.2975388970072714, 0.04244215115672034, 0, 0.1854661576209224, 0.03190823855041458, 0, 0, 0, 0, 0, 0.9720058023732747),
  @(5.155259563753589, 1.423289536696075, 0.09546712611732744, 2.017600804443191, 0.3892102654806422, 0.2568431175060297, 0.03234465346196602, 0, 0.1753424130386563, 0.04039756437714281, 0, 0, 0, 0, 0, 0.8640012070944749),
  @(4.245627476704612, 1.1799504783271, 0.07898959012827333, 1.660708142630455, 0.3321353404837453, 0.2169474456511011, 0.02266980238579031, 0, 0.1666808661305836, 0.05248172212341196, 0, 0, 0, 0, 0, 0.7606831393047884),
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
  $rowIdx = $i + 2
  $vals = $rowsData[$i]
  for ($j = 0; $j -lt $cols.Count; $j++) {
    $ws.Range("$($cols[$j])$rowIdx").Value2 = $vals[$j]
  }
}
